$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 138
$ws.Range("G2").Value = 0.04575443874290686

$ws.Range("E3").Value = 207
$ws.Range("G3").Value = 0.06324657660138239

$ws.Range("E4").Value = 69
$ws.Range("G4").Value = 0.02267945347652807

$ws.Range("E5").Value = 140
$ws.Range("G5").Value = 0.02162823772088771

$ws.Range("E6").Value = 261
$ws.Range("G6").Value = 0.03922053475144366

$ws.Range("E7").Value = 103
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.01061226102384237

$ws.Range("E8").Value = 231
$ws.Range("G8").Value = 0.01810323634571336

$ws.Range("E9").Value = 247
$ws.Range("G9").Value = 0.02660723692926331

$ws.Range("E10").Value = 141
$ws.Range("G10").Value = 0.01378842105240141

$ws.Range("E11").Value = 116
$ws.Range("G11").Value = 0.009469554293510648

$ws.Range("E12").Value = 111
$ws.Range("G12").Value = 0.007281291733332003

$ws.Range("E13").Value = 277
$ws.Range("G13").Value = 0.008819291968009482

$ws.Range("E14").Value = 251
$ws.Range("G14").Value = 0.00598171552042424

$ws.Range("E15").Value = 214
$ws.Range("G15").Value = 0.00670081302902699

$ws.Range("E16").Value = 326
$ws.Range("G16").Value = 0.007949204903071607

$ws.Range("E17").Value = 170
$ws.Range("G17").Value = 0.02598671790384382

$ws.Range("E18").Value = 233
$ws.Range("G18").Value = 0.03692968282165096

$ws.Range("E19").Value = 119
$ws.Range("G19").Value = 0.01366432290567512

$ws.Range("E20").Value = 188
$ws.Range("G20").Value = 0.01301287560234095

$ws.Range("E21").Value = 106

$ws.Range("E22").Value = 176
$ws.Range("G22").Value = 0.01209398928722977
